# Applies the "removendo requisitos que envolvem exames" edit:
#  1) Trim the "a requisição de exames" mention (and the matching "e fazendo os
#     exames necessários..." clause) from the intro paragraph.
#  2) Trim "exames feitos, " from the doctor's-consultation paragraph.
#  3) Remove the "Falta de controle sobre quais exames..." bullet entirely.
#  4) Remove the "cadastrar exames" bullet entirely.
#  5) Remove the "requisitar exames para o paciente" bullet entirely.
#  6) Remove the RF08 (Cadastro de exames) and RF09 (Requisição de exames)
#     rows from the functional-requirements table.

$d = $word.ActiveDocument

# --- 1) Intro paragraph (section 1) -----------------------------------
# Drop "a requisição de exames, " from the list of examples ...
$d.Content.Find.Execute("a requisição de exames, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
# ... and drop the matching "... e fazendo os exames necessários ..." clause.
$d.Content.Find.Execute(" e fazendo os exames necessários para que o profissional da saúde chegue a um diagnóstico mais preciso", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- 2) "Cada médico poderá realizar..." paragraph (section 2) -------
$d.Content.Find.Execute("exames feitos, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- Helper: delete the first whole paragraph whose text contains $marker
function Remove-ParagraphByMarker($marker) {
    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$marker*") {
            $p.Range.Delete()
            return $true
        }
    }
    return $false
}

# --- 3) Remove "Falta de controle sobre quais exames..." bullet -------
Remove-ParagraphByMarker("Falta de controle sobre quais exames") | Out-Null

# --- 4) Remove "Deverá existir uma área para cadastrar exames." bullet -
Remove-ParagraphByMarker("para cadastrar exames") | Out-Null

# --- 5) Remove "...requisitar exames para o paciente." bullet ---------
Remove-ParagraphByMarker("requisitar exames para o paciente") | Out-Null

# --- 6) Remove RF08 and RF09 rows from the requirements table ---------
$tbl = $d.Tables.Item(1)
for ($i = $tbl.Rows.Count; $i -ge 1; $i--) {
    $row = $tbl.Rows.Item($i)
    $label = $row.Cells.Item(1).Range.Text
    if (($label -like "*RF08*") -or ($label -like "*RF09*")) {
        $row.Delete()
    }
}
